$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Glicemia"
$ws.Range("D1").Value = "Urea"
$ws.Range("E1").Value = "ALT"
$ws.Range("F1").Value = "AST"

# --- Data rows (ID, Type, Glicemia, Urea, ALT, AST) ---
$data = @(
    @(124, 0, 78, 25, 7, $null),
    @(132, 0, 65, 5, 24, 8),
    @(154, 0, 115, 11, 63, 16),
    @(163, 1, 70, 21, 23, 13),
    @(174, 1, 46, 16, 43, 6),
    @(176, 0, $null, 13, 13, 29),
    @(234, 0, 86, 4, 53, 33),
    @(241, 0, 89, 5, 67, 23),
    @(246, 1, 70, 9, 32, 5),
    @(254, 1, 69, 14, 14, 23),
    @(261, 0, 48, 19, 12, 15),
    @(261, 0, 48, 19, 12, 15),
    @(275, 0, 78, 6, 13, 64),
    @(290, 0, 146, 35, 42, 52),
    @(347, 1, 67, 8, 25, 23),
    @(423, 1, 95, 8, 64, 70),
    @(467, 1, 135, 10, 7, 8),
    @(567, 1, 75, $null, $null, $null),
    @(728, 1, 100, 31, 5, 32),
    @(735, 0, 72, 52, 6, 25),
    @(824, 1, 99, 8, 49, 10),
    @(865, 1, 79, 12, 38, 11),
    @(920, 0, 86, 15, 17, 20),
    @(124, 0, 78, 25, 29, 30)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    if ($r[2] -ne $null) { $ws.Cells.Item($row, 3).Value = $r[2] } else { $ws.Cells.Item($row, 3).ClearContents() }
    if ($r[3] -ne $null) { $ws.Cells.Item($row, 4).Value = $r[3] } else { $ws.Cells.Item($row, 4).ClearContents() }
    if ($r[4] -ne $null) { $ws.Cells.Item($row, 5).Value = $r[4] } else { $ws.Cells.Item($row, 5).ClearContents() }
    if ($r[5] -ne $null) { $ws.Cells.Item($row, 6).Value = $r[5] } else { $ws.Cells.Item($row, 6).ClearContents() }
    $row++
}

# --- Window / view state ---
try {
    $win = $excel.ActiveWindow
    $win.Left = 57480
    $win.Top = 8385
    $win.Width = 29040
    $win.Height = 17520
} catch {
    # Window geometry not settable in this headless runtime; ignore.
}

# --- Selection ---
$ws.Range("K23").Select()

$wb.Save()
